# Shift the whole historical data set forward by one calendar day.
#
# Column A (Timestamp, CET) holds date/time serials -> add 1 (day) to each.
# Column E (Lookup) holds strings like "25.01.20261" (DD.MM.YYYY + running
# quarter-hour index within the day) -> bump the DD.MM.YYYY part by one day
# while leaving the trailing index untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value2 = $aCell.Value2 + 1

    $eCell = $ws.Cells.Item($r, 5)
    $lookup = [string]$eCell.Value2
    if ($lookup -match '^(\d{2})\.(\d{2})\.(\d{4})(\d+)$') {
        $day = [int]$matches[1]
        $month = [int]$matches[2]
        $year = [int]$matches[3]
        $suffix = $matches[4]

        $dt = Get-Date -Year $year -Month $month -Day $day
        $dt = $dt.AddDays(1)

        $eCell.Value2 = $dt.ToString("dd.MM.yyyy") + $suffix
    }
}
